$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'282.28"
$ws.Range("D3").Value = "'20.59"
$ws.Range("D4").Value = "'6.215"
$ws.Range("D5").Value = "'0.06170"
$ws.Range("D6").Value = "'3.581"
$ws.Range("D7").Value = "'1.513"
$ws.Range("D8").Value = "'6.551"
$ws.Range("D9").Value = "'0.8185"
$ws.Range("D10").Value = "'0.01385"
$ws.Range("D11").Value = "'0.1640"
$ws.Range("D12").Value = "'0.08433"
$ws.Range("D13").Value = "'0.03524"
$ws.Range("D14").Value = "'0.03216"
$ws.Range("D15").Value = "'0.09151"
$ws.Range("D16").Value = "'3.714"
$ws.Range("D17").Value = "'0.001646"
$ws.Range("D18").Value = "'0.04716"
$ws.Range("D19").Value = "'0.006448"
$ws.Range("D20").Value = "'0.006170"
$ws.Range("D23").Value = "'3.831"
$ws.Range("D24").Value = "'2.343"
$ws.Range("D25").Value = "'0.3314"
$ws.Range("D26").Value = "'0.1251"
$ws.Range("D40").Value = "'0.04708"
$ws.Range("D41").Value = "'0.007207"
$ws.Range("D43").Value = "'0.1098"
$ws.Range("D44").Value = "'0.01138"
$ws.Range("D45").Value = "'0.00006609"
$ws.Range("D47").Value = "'0.8011"
$ws.Range("D48").Value = "'0.002859"
